$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Teachers")

$ws.Cells.Item(8, 1).Value = "samikb@gmail.com"
$ws.Cells.Item(8, 2).Value = "ABCD"
$ws.Cells.Item(8, 3).Value = '$2b$12$Mbu0Y4wabzc0HchfxKPOIukURdO71QwM2sHmt4U.esjDlNzJkRyi2'
$ws.Cells.Item(8, 4).Value = "teacher"
$ws.Cells.Item(8, 5).Value = "DSGT, DLCOA"
